$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '65.916.94'
Set-TextValue 'E2' '  -3.39%  '
Set-TextValue 'D3' '2.403.74'
Set-TextValue 'E3' '  -5.10%  '
Set-TextValue 'E4' '  +0.27%  '
Set-TextValue 'D5' '567.79'
Set-TextValue 'E5' '  -4.42%  '
Set-TextValue 'D6' '160.54'
Set-TextValue 'E6' '  -9.68%  '
Set-TextValue 'E7' '  +0.29%  '
Set-TextValue 'D8' '0.497'
Set-TextValue 'E8' '  -6.53%  '
Set-TextValue 'D9' '2.400.86'
Set-TextValue 'E9' '  -5.19%  '
Set-TextValue 'D10' '0.130'
Set-TextValue 'E10' '  -8.98%  '
Set-TextValue 'D11' '0.163'
Set-TextValue 'E11' '  -0.87%  '
Set-TextValue 'D12' '0.324'
Set-TextValue 'E12' '  -6.57%  '
Set-TextValue 'D13' '4.72'
Set-TextValue 'E13' '  -7.86%  '
Set-TextValue 'D14' '24.69'
Set-TextValue 'E14' '  -8.19%  '
Set-TextValue 'B15' 'WrappedBTC'
Set-TextValue 'C15' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D15' '65.976.01'
Set-TextValue 'E15' '  -3.29%  '
Set-TextValue 'B16' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C16' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D16' '2.742.05'
Set-TextValue 'E16' '  -8.40%  '
Set-TextValue 'B17' 'ShibaInu'
Set-TextValue 'C17' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D17' '0.0000163'
Set-TextValue 'E17' '  -9.51%  '
Set-TextValue 'D18' '2.430.06'
Set-TextValue 'E18' '  -4.13%  '
Set-TextValue 'D19' '11.03'
Set-TextValue 'E19' '  -4.55%  '
Set-TextValue 'D20' '7.42'
Set-TextValue 'E20' '  -7.13%  '
Set-TextValue 'D21' '347.75'
Set-TextValue 'E21' '  -5.56%  '
Set-TextValue 'D22' '3.93'
Set-TextValue 'E22' '  -6.50%  '
Set-TextValue 'D23' '1.00'
Set-TextValue 'E23' '  +0.05%  '
Set-TextValue 'B24' 'Litecoin'
Set-TextValue 'C24' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D24' '69.27'
Set-TextValue 'E24' '  -2.30%  '
Set-TextValue 'B25' 'NEARProtocol'
Set-TextValue 'C25' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D25' '4.13'
Set-TextValue 'E25' '  -12.40%  '
Set-TextValue 'D26' '1.72'
Set-TextValue 'E26' '  -11.25%  '
Set-TextValue 'D27' '8.80'
Set-TextValue 'E27' '  -13.90%  '
Set-TextValue 'D28' '0.994'
Set-TextValue 'E28' '  -0.34%  '
Set-TextValue 'D29' '2.588.42'
Set-TextValue 'E29' '  -2.72%  '
Set-TextValue 'D30' '0.0₃0873'
Set-TextValue 'E30' '  -12.50%  '
Set-TextValue 'D31' '7.60'
Set-TextValue 'E31' '  -8.30%  '
Set-TextValue 'D32' '475.18'
Set-TextValue 'E32' '  -12.28%  '
Set-TextValue 'D33' '1.76'
Set-TextValue 'E33' '  -6.43%  '
Set-TextValue 'B34' 'FirstDigitalUSD'
Set-TextValue 'C34' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D34' '1.00'
Set-TextValue 'E34' '  +0.44%  '
Set-TextValue 'B35' 'Fetch.AI'
Set-TextValue 'C35' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D35' '1.15'
Set-TextValue 'E35' '  -14.00%  '
Set-TextValue 'B36' 'Monero'
Set-TextValue 'C36' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D36' '157.64'
Set-TextValue 'E36' '  +0.33%  '
Set-TextValue 'B37' 'Kaspa'
Set-TextValue 'C37' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D37' '0.115'
Set-TextValue 'E37' '  -11.22%  '
Set-TextValue 'B38' 'WhiteBITCoin'
Set-TextValue 'C38' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D38' '18.52'
Set-TextValue 'E38' '  -0.96%  '
Set-TextValue 'B39' 'EthereumClassic'
Set-TextValue 'C39' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D39' '18.17'
Set-TextValue 'E39' '  -3.68%  '
Set-TextValue 'B40' 'ImmutableX'
Set-TextValue 'C40' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D40' '1.32'
Set-TextValue 'E40' '  -9.73%  '
Set-TextValue 'B41' 'PolygonEcosystemToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D41' '0.313'
Set-TextValue 'E41' '  -12.10%  '
Set-TextValue 'B42' 'RenderToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue 'D42' '4.58'
Set-TextValue 'E42' '  -12.24%  '
Set-TextValue 'B43' 'Stacks'
Set-TextValue 'C43' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D43' '1.58'
Set-TextValue 'E43' '  -12.59%  '
Set-TextValue 'B44' 'OKB'
Set-TextValue 'C44' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D44' '39.01'
Set-TextValue 'E44' '  -2.55%  '
Set-TextValue 'D45' '2.27'
Set-TextValue 'E45' '  -11.44%  '
Set-TextValue 'B46' 'Aave'
Set-TextValue 'C46' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D46' '135.72'
Set-TextValue 'E46' '  -7.87%  '
Set-TextValue 'B47' 'Filecoin'
Set-TextValue 'C47' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D47' '3.41'
Set-TextValue 'E47' '  -8.58%  '
Set-TextValue 'B48' 'ARBITRUM'
Set-TextValue 'C48' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D48' '0.502'
Set-TextValue 'E48' '  -10.63%  '
Set-TextValue 'B49' 'BabyDogeCoin'
Set-TextValue 'C49' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D49' '0.0₆0249'
Set-TextValue 'E49' '  -10.40%  '
Set-TextValue 'B50' 'Cronos'
Set-TextValue 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D50' '0.0721'
Set-TextValue 'E50' '  -4.71%  '
Set-TextValue 'D51' '1.54'
Set-TextValue 'E51' '  -9.39%  '
